$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a piece of text into a cell WITHOUT letting Excel's
# value-assignment auto-detect it as a number/date and reassign the
# cell's style. A value-only paste (xlPasteValues) from a scratch cell
# leaves the destination's existing style (quote-prefix / number format /
# font) completely untouched, which is what plain `.Value = ...` does not
# do for strings that look like times (e.g. "09:00 AM").
function Set-TextKeepStyle($sheet, $targetAddr, $text) {
    $scratch = $sheet.Range("D5")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $sheet.Range($targetAddr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Insert a new column before column C (shifts C:M -> D:N), carrying over
# formatting/widths from the original column C, just like Excel's
# "Insert Sheet Columns" command.
$ws.Columns("C").EntireColumn.Insert()

# New column header + values ("AddressLine2" / "address line 2")
$ws.Range("C1").Value = "AddressLine2"
$ws.Range("C2").Value = "address line 2"
$ws.Range("C3").Value = "address line 2"
$ws.Range("C4").Value = "address line 2"

# Column B ("AddressLine1" data) text tweak: "address line 1" -> "address line1"
$ws.Range("B2").Value = "address line1"
$ws.Range("B3").Value = "address line1"
$ws.Range("B4").Value = "address line1"

# Updated pickup/drop-off times for rows 2 and 3 (columns shifted right by
# the column insert: old E:H is now F:I). These are text-formatted cells
# (quote-prefixed) in the original, so write them via the style-preserving
# helper rather than a plain .Value assignment.
Set-TextKeepStyle $ws "F2" "09:00 AM"
Set-TextKeepStyle $ws "G2" "04:00 AM"
Set-TextKeepStyle $ws "H2" "09:00 AM"
Set-TextKeepStyle $ws "I2" "04:00 AM"

Set-TextKeepStyle $ws "F3" "09:00 AM"
Set-TextKeepStyle $ws "G3" "04:00 AM"
Set-TextKeepStyle $ws "H3" "09:00 AM"
Set-TextKeepStyle $ws "I3" "04:00 AM"

# Update the active selection to D2, matching the saved view state
$ws.Range("D2").Select()
